# ------------------------------------------------------------------
# Registra um novo "registro" de trabalho (Elétrica / VDS) para o ID 010
# e atualiza os resumos derivados (Gráficos) e o orçamento (Orçamentos).
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ====================================================================
# 1) Sheet "Registros": append new row 7 with the new time entry
# ====================================================================
$wsReg = $wb.Worksheets.Item("Registros")
$r = 7

# Text-looking values that Excel would otherwise auto-coerce (dates /
# numbers with leading zeros) must be forced to Text first so the
# literal string is preserved, exactly like the other rows in the sheet.
$wsReg.Cells.Item($r, 1).NumberFormat = "@"
$wsReg.Cells.Item($r, 1).Value = "2025-10-27"

$wsReg.Cells.Item($r, 2).NumberFormat = "@"
$wsReg.Cells.Item($r, 2).Value = "010"

$wsReg.Cells.Item($r, 3).Value = "JOSE GENILSOS MARTINS SOARES"
$wsReg.Cells.Item($r, 4).Value = "Elétrica"
$wsReg.Cells.Item($r, 5).Value = "VDS"

$wsReg.Cells.Item($r, 6).NumberFormat = "@"
$wsReg.Cells.Item($r, 6).Value = "15"

$wsReg.Cells.Item($r, 7).Value = "15:40:00"
$wsReg.Cells.Item($r, 8).Value = "16:40:00"
$wsReg.Cells.Item($r, 9).Value = "registro"

# ====================================================================
# 2) Sheet "Orçamentos": append new budget row 4 (Hidráulica / LOGISTÍCA)
# ====================================================================
$wsOrc = $wb.Worksheets.Item("Orçamentos")

$wsOrc.Cells.Item(4, 1).Value = "Hidráulica"
$wsOrc.Cells.Item(4, 2).Value = "LOGISTÍCA"

$wsOrc.Cells.Item(4, 3).NumberFormat = "@"
$wsOrc.Cells.Item(4, 3).Value = "1"

$wsOrc.Cells.Item(4, 4).NumberFormat = "0"
$wsOrc.Cells.Item(4, 4).Value = 20

# ====================================================================
# 3) Sheet "Gráficos": re-sort the summary (alphabetical) and recompute
#    the derived hours for the affected projects.
#    Final order: Hidráulica - LOGISTÍCA - 1 (row2, new)
#                 Marcenaria Estrutural - CDS - 12 (row3, updated)
#                 Marcenaria Móvel - ODS - 12 (row4, unchanged, moved)
# ====================================================================
$wsGra = $wb.Worksheets.Item("Gráficos")

# 3a) Move current row 3 (Marcenaria Móvel) down to row 4 first so it is
#     not clobbered when row 2/3 are rewritten below.
$wsGra.Cells.Item(4, 1).Value = "Marcenaria Móvel - ODS - 12"
$wsGra.Cells.Item(4, 2).NumberFormat = "0.00"
$wsGra.Cells.Item(4, 2).Value = 0
$wsGra.Cells.Item(4, 3).NumberFormat = "0"
$wsGra.Cells.Item(4, 3).Value = 23
$wsGra.Cells.Item(4, 4).NumberFormat = "0.00"
$wsGra.Cells.Item(4, 4).Value = 23
$wsGra.Cells.Item(4, 5).NumberFormat = "0.00"
$wsGra.Cells.Item(4, 5).Value = 0

# 3b) Row 3: Marcenaria Estrutural - CDS - 12, now with 1h worked
$wsGra.Cells.Item(3, 1).Value = "Marcenaria Estrutural - CDS - 12"
$wsGra.Cells.Item(3, 2).NumberFormat = "0.00"
$wsGra.Cells.Item(3, 2).Value = 1
$wsGra.Cells.Item(3, 3).NumberFormat = "0"
$wsGra.Cells.Item(3, 3).Value = 30
$wsGra.Cells.Item(3, 4).NumberFormat = "0.00"
$wsGra.Cells.Item(3, 4).Value = 29
$wsGra.Cells.Item(3, 5).NumberFormat = "0.00"
$wsGra.Cells.Item(3, 5).Value = 3.33

# 3c) Row 2: Hidráulica - LOGISTÍCA - 1 (new project, nothing worked yet)
$wsGra.Cells.Item(2, 1).Value = "Hidráulica - LOGISTÍCA - 1"
$wsGra.Cells.Item(2, 2).NumberFormat = "0.00"
$wsGra.Cells.Item(2, 2).Value = 0
$wsGra.Cells.Item(2, 3).NumberFormat = "0"
$wsGra.Cells.Item(2, 3).Value = 20
$wsGra.Cells.Item(2, 4).NumberFormat = "0.00"
$wsGra.Cells.Item(2, 4).Value = 20
$wsGra.Cells.Item(2, 5).NumberFormat = "0.00"
$wsGra.Cells.Item(2, 5).Value = 0

# ====================================================================
# 4) Extend conditional formatting to cover the new row (E2:E4/B2:B4/C2:C4)
# ====================================================================
$null = $wsGra.Range("E2:E3").FormatConditions.AddColorScale(3)
$null = $wsGra.Range("B2:B3").FormatConditions.AddDatabar()
$null = $wsGra.Range("C2:C3").FormatConditions.AddDatabar()
$null = $wsGra.Range("E2:E4").FormatConditions.AddColorScale(3)
$null = $wsGra.Range("B2:B4").FormatConditions.AddDatabar()
$null = $wsGra.Range("C2:C4").FormatConditions.AddDatabar()

# ====================================================================
# 5) Update the two charts so their series cover rows 2:4 instead of 2:3
# ====================================================================
$chart1 = $wsGra.ChartObjects().Item(1).Chart
$chart1.SeriesCollection().Item(1).Formula = "=SERIES(,'Gráficos'!`$A`$2:`$A`$4,'Gráficos'!`$B`$2:`$B`$4,1)"
$chart1.SeriesCollection().Item(2).Formula = "=SERIES(,'Gráficos'!`$A`$2:`$A`$4,'Gráficos'!`$C`$2:`$C`$4,2)"

$chart2 = $wsGra.ChartObjects().Item(2).Chart
$chart2.SeriesCollection().Item(1).Formula = "=SERIES(,'Gráficos'!`$A`$2:`$A`$4,'Gráficos'!`$E`$2:`$E`$4,1)"
$chart2.SeriesCollection().Item(2).Formula = "=SERIES(,'Gráficos'!`$A`$2:`$A`$4,'Gráficos'!`$D`$2:`$D`$4,2)"
